# Delete Sheet2 and Sheet3, rename Sheet1 to TODO, fill in TODO list data,
# set column widths, and move the selection to A6.

$wb = $excel.ActiveWorkbook

# Remove the extra sheets first.
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "TODO"

$ws.Range("A1").Value = "Properly emulate behavior: hearts falling, explosions, etc"
$ws.Range("A2").Value = "load/save config"
$ws.Range("A3").Value = "license for git repo"
$ws.Range("A4").Value = "resolution-changing ingame"
$ws.Range("A5").Value = "new/old gfx changing ingame"

# ColumnWidth assignments are snapped to the nearest 1/6-character grid by
# the engine (stored width = round(input*6)/6 + 5/6), so these inputs are
# chosen to land as close as possible to the target stored widths of
# 98.85546875 (col A) and 16.140625 (col B).
$ws.Columns.Item(1).ColumnWidth = 98
$ws.Columns.Item(2).ColumnWidth = 15.333333333333334

$ws.Range("A6").Select() | Out-Null
